# Apply scheduled-runner price/profit refresh to the Leve profit sheets.
# Each (sheet, cell) pair below is set to the newly-fetched market value.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, "H").Value = 364.69513
$ws.Cells.Item(17, "I").Value = 700
$ws.Cells.Item(17, "J").Value = 356.3125
$ws.Cells.Item(17, "K").Value = 2100
$ws.Cells.Item(17, "L").Value = 1068.9375
$ws.Cells.Item(17, "M").Value = -1932
$ws.Cells.Item(17, "N").Value = -1404.9375
$ws.Cells.Item(40, "H").Value = 1872
$ws.Cells.Item(40, "I").Value = 3500
$ws.Cells.Item(40, "J").Value = 1465
$ws.Cells.Item(40, "K").Value = 3500
$ws.Cells.Item(40, "L").Value = 1465
$ws.Cells.Item(40, "M").Value = -3325
$ws.Cells.Item(40, "N").Value = -1815
$ws.Cells.Item(64, "H").Value = 4170.684
$ws.Cells.Item(64, "I").Value = 4308
$ws.Cells.Item(64, "J").Value = 4070.818
$ws.Cells.Item(64, "K").Value = 4308
$ws.Cells.Item(64, "L").Value = 4070.818
$ws.Cells.Item(64, "M").Value = -4060
$ws.Cells.Item(64, "N").Value = -4566.818
$ws.Cells.Item(67, "H").Value = 4170.684
$ws.Cells.Item(67, "I").Value = 4308
$ws.Cells.Item(67, "J").Value = 4070.818
$ws.Cells.Item(67, "K").Value = 4308
$ws.Cells.Item(67, "L").Value = 4070.818
$ws.Cells.Item(67, "M").Value = -3450
$ws.Cells.Item(67, "N").Value = -5786.818
$ws.Cells.Item(129, "H").Value = 597.2857
$ws.Cells.Item(129, "I").Value = 386.2
$ws.Cells.Item(129, "J").Value = 1125
$ws.Cells.Item(129, "K").Value = 1158.6
$ws.Cells.Item(129, "L").Value = 3375
$ws.Cells.Item(129, "M").Value = 3841.4
$ws.Cells.Item(129, "N").Value = -13375
$ws.Cells.Item(137, "H").Value = 235237.88
$ws.Cells.Item(137, "I").Value = 298938.44
$ws.Cells.Item(137, "J").Value = 39236.08
$ws.Cells.Item(137, "K").Value = 896815.3200000001
$ws.Cells.Item(137, "L").Value = 117708.24
$ws.Cells.Item(137, "M").Value = -894265.3200000001
$ws.Cells.Item(137, "N").Value = -122808.24
$ws.Cells.Item(140, "H").Value = 20000
$ws.Cells.Item(140, "J").Value = 20000
$ws.Cells.Item(140, "L").Value = 20000
$ws.Cells.Item(140, "N").Value = -30360

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, "H").Value = 548654.5
$ws.Cells.Item(32, "I").Value = 2838.7144
$ws.Cells.Item(32, "K").Value = 2838.7144
$ws.Cells.Item(32, "M").Value = -2551.7144
$ws.Cells.Item(61, "H").Value = 6383.074
$ws.Cells.Item(61, "I").Value = 6879.273
$ws.Cells.Item(61, "J").Value = 4199.8
$ws.Cells.Item(61, "K").Value = 6879.273
$ws.Cells.Item(61, "L").Value = 4199.8
$ws.Cells.Item(61, "M").Value = -6667.273
$ws.Cells.Item(61, "N").Value = -4623.8
$ws.Cells.Item(74, "H").Value = 3347.8
$ws.Cells.Item(74, "I").Value = 579.6786
$ws.Cells.Item(74, "J").Value = 9806.75
$ws.Cells.Item(74, "K").Value = 579.6786
$ws.Cells.Item(74, "L").Value = 9806.75
$ws.Cells.Item(74, "M").Value = 294.3214
$ws.Cells.Item(74, "N").Value = -11554.75
$ws.Cells.Item(77, "H").Value = 3347.8
$ws.Cells.Item(77, "I").Value = 579.6786
$ws.Cells.Item(77, "J").Value = 9806.75
$ws.Cells.Item(77, "K").Value = 2898.393
$ws.Cells.Item(77, "L").Value = 49033.75
$ws.Cells.Item(77, "M").Value = 1469.607
$ws.Cells.Item(77, "N").Value = -57769.75
$ws.Cells.Item(95, "H").Value = 17400
$ws.Cells.Item(95, "J").Value = 17400
$ws.Cells.Item(95, "L").Value = 17400
$ws.Cells.Item(95, "N").Value = -22892
$ws.Cells.Item(132, "H").Value = 2718736.2
$ws.Cells.Item(132, "I").Value = 3290312.8
$ws.Cells.Item(132, "J").Value = 3748.25
$ws.Cells.Item(132, "K").Value = 9870938.399999999
$ws.Cells.Item(132, "L").Value = 11244.75
$ws.Cells.Item(132, "M").Value = -9868408.399999999
$ws.Cells.Item(132, "N").Value = -16304.75
$ws.Cells.Item(136, "H").Value = 6383.074
$ws.Cells.Item(136, "I").Value = 6879.273
$ws.Cells.Item(136, "J").Value = 4199.8
$ws.Cells.Item(136, "K").Value = 20637.819
$ws.Cells.Item(136, "L").Value = 12599.4
$ws.Cells.Item(136, "M").Value = -18087.819
$ws.Cells.Item(136, "N").Value = -17699.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, "H").Value = 2209.0908
$ws.Cells.Item(105, "I").Value = 2214.2856
$ws.Cells.Item(105, "J").Value = 2200
$ws.Cells.Item(105, "K").Value = 2214.2856
$ws.Cells.Item(105, "L").Value = 2200
$ws.Cells.Item(105, "M").Value = -467.2856000000002
$ws.Cells.Item(105, "N").Value = -5694
$ws.Cells.Item(134, "H").Value = 8131581.5
$ws.Cells.Item(134, "I").Value = 10417947
$ws.Cells.Item(134, "J").Value = 2279.3333
$ws.Cells.Item(134, "K").Value = 31253841
$ws.Cells.Item(134, "L").Value = 6837.999899999999
$ws.Cells.Item(134, "M").Value = -31251306
$ws.Cells.Item(134, "N").Value = -11907.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, "H").Value = 11018
$ws.Cells.Item(31, "I").Value = 5524.773
$ws.Cells.Item(31, "J").Value = 24445.889
$ws.Cells.Item(31, "K").Value = 5524.773
$ws.Cells.Item(31, "L").Value = 24445.889
$ws.Cells.Item(31, "M").Value = -5229.773
$ws.Cells.Item(31, "N").Value = -25035.889
$ws.Cells.Item(34, "H").Value = 11018
$ws.Cells.Item(34, "I").Value = 5524.773
$ws.Cells.Item(34, "J").Value = 24445.889
$ws.Cells.Item(34, "K").Value = 5524.773
$ws.Cells.Item(34, "L").Value = 24445.889
$ws.Cells.Item(34, "M").Value = -5322.773
$ws.Cells.Item(34, "N").Value = -24849.889
$ws.Cells.Item(58, "H").Value = 4934109.5
$ws.Cells.Item(58, "I").Value = 7143819.5
$ws.Cells.Item(58, "J").Value = 23641.666
$ws.Cells.Item(58, "K").Value = 7143819.5
$ws.Cells.Item(58, "L").Value = 23641.666
$ws.Cells.Item(58, "M").Value = -7143616.5
$ws.Cells.Item(58, "N").Value = -24047.666
$ws.Cells.Item(107, "H").Value = 489.72726
$ws.Cells.Item(107, "I").Value = 435.16666
$ws.Cells.Item(107, "J").Value = 555.2
$ws.Cells.Item(107, "K").Value = 435.16666
$ws.Cells.Item(107, "L").Value = 555.2
$ws.Cells.Item(107, "M").Value = 1484.83334
$ws.Cells.Item(107, "N").Value = -4395.2
$ws.Cells.Item(132, "H").Value = 11910059
$ws.Cells.Item(132, "I").Value = 19608744
$ws.Cells.Item(132, "J").Value = 12092
$ws.Cells.Item(132, "K").Value = 58826232
$ws.Cells.Item(132, "L").Value = 36276
$ws.Cells.Item(132, "M").Value = -58823702
$ws.Cells.Item(132, "N").Value = -41336
$ws.Cells.Item(134, "H").Value = 8447670
$ws.Cells.Item(134, "I").Value = 8622405
$ws.Cells.Item(134, "J").Value = 7814253.5
$ws.Cells.Item(134, "K").Value = 25867215
$ws.Cells.Item(134, "L").Value = 23442760.5
$ws.Cells.Item(134, "M").Value = -25864680
$ws.Cells.Item(134, "N").Value = -23447830.5
$ws.Cells.Item(136, "H").Value = 4934109.5
$ws.Cells.Item(136, "I").Value = 7143819.5
$ws.Cells.Item(136, "J").Value = 23641.666
$ws.Cells.Item(136, "K").Value = 21431458.5
$ws.Cells.Item(136, "L").Value = 70924.99800000001
$ws.Cells.Item(136, "M").Value = -21428908.5
$ws.Cells.Item(136, "N").Value = -76024.99800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(40, "H").Value = 125
$ws.Cells.Item(40, "J").Value = 150
$ws.Cells.Item(40, "L").Value = 600
$ws.Cells.Item(40, "N").Value = -738
$ws.Cells.Item(69, "H").Value = 1844.7142
$ws.Cells.Item(69, "I").Value = 399
$ws.Cells.Item(69, "J").Value = 2085.6667
$ws.Cells.Item(69, "K").Value = 1197
$ws.Cells.Item(69, "L").Value = 6257.000100000001
$ws.Cells.Item(69, "M").Value = -386
$ws.Cells.Item(69, "N").Value = -7879.000100000001
$ws.Cells.Item(72, "H").Value = 1844.7142
$ws.Cells.Item(72, "I").Value = 399
$ws.Cells.Item(72, "J").Value = 2085.6667
$ws.Cells.Item(72, "K").Value = 3591
$ws.Cells.Item(72, "L").Value = 18771.0003
$ws.Cells.Item(72, "M").Value = 465
$ws.Cells.Item(72, "N").Value = -26883.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, "H").Value = 74076040
$ws.Cells.Item(132, "I").Value = 111111770
$ws.Cells.Item(132, "J").Value = 4578.8887
$ws.Cells.Item(132, "K").Value = 333335310
$ws.Cells.Item(132, "L").Value = 13736.6661
$ws.Cells.Item(132, "M").Value = -333332780
$ws.Cells.Item(132, "N").Value = -18796.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, "H").Value = 3892.7058
$ws.Cells.Item(82, "I").Value = 1325
$ws.Cells.Item(82, "J").Value = 4682.769
$ws.Cells.Item(82, "K").Value = 1325
$ws.Cells.Item(82, "L").Value = 4682.769
$ws.Cells.Item(82, "M").Value = -964
$ws.Cells.Item(82, "N").Value = -5404.769
$ws.Cells.Item(85, "H").Value = 3892.7058
$ws.Cells.Item(85, "I").Value = 1325
$ws.Cells.Item(85, "J").Value = 4682.769
$ws.Cells.Item(85, "K").Value = 1325
$ws.Cells.Item(85, "L").Value = 4682.769
$ws.Cells.Item(85, "M").Value = -77
$ws.Cells.Item(85, "N").Value = -7178.769
$ws.Cells.Item(100, "H").Value = 2025
$ws.Cells.Item(100, "I").Value = 2116.6667
$ws.Cells.Item(100, "K").Value = 2116.6667
$ws.Cells.Item(100, "M").Value = -1575.6667
$ws.Cells.Item(132, "H").Value = 5884560
$ws.Cells.Item(132, "I").Value = 11112833
$ws.Cells.Item(132, "J").Value = 2753
$ws.Cells.Item(132, "K").Value = 33338499
$ws.Cells.Item(132, "L").Value = 8259
$ws.Cells.Item(132, "M").Value = -33335969
$ws.Cells.Item(132, "N").Value = -13319
$ws.Cells.Item(136, "H").Value = 3031.8572
$ws.Cells.Item(136, "I").Value = 3256.9473
$ws.Cells.Item(136, "J").Value = 2556.6667
$ws.Cells.Item(136, "K").Value = 9770.841899999999
$ws.Cells.Item(136, "L").Value = 7670.000100000001
$ws.Cells.Item(136, "M").Value = -7220.841899999999
$ws.Cells.Item(136, "N").Value = -12770.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, "H").Value = 43951836
$ws.Cells.Item(132, "I").Value = 44446396
$ws.Cells.Item(132, "J").Value = 42839070
$ws.Cells.Item(132, "K").Value = 133339188
$ws.Cells.Item(132, "L").Value = 128517210
$ws.Cells.Item(132, "M").Value = -133336658
$ws.Cells.Item(132, "N").Value = -128522270
$ws.Cells.Item(136, "H").Value = 29524916
$ws.Cells.Item(136, "I").Value = 19788592
$ws.Cells.Item(136, "J").Value = 45457084
$ws.Cells.Item(136, "K").Value = 59365776
$ws.Cells.Item(136, "L").Value = 136371252
$ws.Cells.Item(136, "M").Value = -59363226
$ws.Cells.Item(136, "N").Value = -136376352
